# Update cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.833.32'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '2.088.97'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.35'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.83'
$ws.Range("E7").Value = '  +2.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.394'
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("E11").Value = '  +2.53%  '

$ws.Range("D12").Value = '2.397.97'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.77'
$ws.Range("E13").Value = '  +2.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.21'
$ws.Range("E14").Value = '  +1.94%  '

$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("E16").Value = '  +1.40%  '

$ws.Range("D17").Value = '2.085.03'
$ws.Range("E17").Value = '  +0.31%  '

$ws.Range("D18").Value = '37.748.40'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.43'

$ws.Range("D21").Value = '0.0₃0838'
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.12'
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  -1.14%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.99'
$ws.Range("E26").Value = '  +0.74%  '

$ws.Range("E27").Value = '  +5.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.03'
$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.59'
$ws.Range("E29").Value = '  +2.10%  '

$ws.Range("E30").Value = '  -1.01%  '

$ws.Range("E31").Value = '  +2.23%  '

$ws.Range("E32").Value = '  +2.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0635'
$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.72'
$ws.Range("E34").Value = '  +3.04%  '

$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("E36").Value = '  +3.04%  '

$ws.Range("E37").Value = '  +3.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("E39").Value = '  -4.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0986'
$ws.Range("E40").Value = '  +3.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.05'
$ws.Range("E41").Value = '  +1.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.92'
$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("E43").Value = '  +1.30%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.464.29'
$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.34'
$ws.Range("E45").Value = '  +1.54%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.07'
$ws.Range("E47").Value = '  +3.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.01'
$ws.Range("E48").Value = '  +3.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.46'
$ws.Range("E49").Value = '  +2.63%  '

$ws.Range("E50").Value = '  +2.11%  '

$ws.Range("D51").Value = '2.281.98'
$ws.Range("E51").Value = '  +0.39%  '
